# TestDataShareSkill.xlsx - "ShareSkill" sheet: align the Startdate/Enddate
# test values on the same day (12082022) instead of two different days
# (12042022 / 12052022).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ShareSkill")

$ws.Range("H2").Value = "12082022"
$ws.Range("I2").Value = "12082022"

# Leave the cursor on the sheet roughly where it ended up in the source
# commit (scrolled right, S2 selected).
$ws.Activate()
$ws.Range("S2").Select() | Out-Null
